# Sync the "Rules" sheet of the non-localizable CodeQuality rules workbook.
# The "BannedPaths" rule row is removed from its old position (row 35) and
# a renamed/updated "BannedPath" rule (no tags, Critical severity) is
# (re)inserted further down the table, just before the "AEM Rules:AEM-3" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Remove the old "BannedPaths" row (row 35); everything below shifts up.
$ws.Rows("35:35").Delete()

# Insert a new row at what is now row 40 (just above "AEM Rules:AEM-3")
# and populate it with the renamed rule. Note: no Tags (column E) value.
$ws.Rows("40:40").Insert()
$ws.Range("A40").Value = "BannedPath"
$ws.Range("B40").Value = "Customer packages should not install content under /libs"
$ws.Range("C40").Value = "Bug"
$ws.Range("D40").Value = "Critical"

# Update the saved selection/active cell on the sheet.
$ws.Activate()
$ws.Range("A37").Select()
